$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (date style for col A, integer style for cols B/C) down
# to the new rows by copying the formats from the last existing row (561).
$ws.Range("A561:C561").Copy()
$ws.Range("A562:C567").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$data = @(
    @(562, 45536, 3.0439546237485131, 3),
    @(563, 45566, 2.8981559020643028, 3),
    @(564, 45597, 2.0403693875787057, 2),
    @(565, 45627, 2.288815139456557, 3),
    @(566, 45658, 3.6692668713944512, 3),
    @(567, 45689, 4.0335035493262668, 3)
)

foreach ($row in $data) {
    $r = $row[0]
    $dateSerial = $row[1]
    $meanVal = $row[2]
    $medVal = $row[3]

    $ws.Cells.Item($r, 1).Value = $dateSerial
    $ws.Cells.Item($r, 2).Value = $meanVal
    $ws.Cells.Item($r, 3).Value = $medVal
}
